$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.975.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.09%  '
$ws.Range("D3").Value = "'1.900.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.03%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.7511"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.67%  '
$ws.Range("D6").Value = "'241.06"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.38%  '
$ws.Range("D7").Value = "'1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '
$ws.Range("D8").Value = "'0.3057"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.35%  '
$ws.Range("D9").Value = "'25.19"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.39%  '
$ws.Range("D10").Value = "'0.06839"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.97%  '
$ws.Range("D11").Value = "'0.07996"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").Value = "'0.7523"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("D13").Value = "'1.899.13"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.74%  '
$ws.Range("D14").Value = "'5.217"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.59%  '
$ws.Range("D15").Value = "'91.20"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = "'6.146"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.14%  '
$ws.Range("D17").Value = "'29.988.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.07%  '
$ws.Range("D18").Value = "'13.94"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.96%  '
$ws.Range("D19").Value = "'0.000007681"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.26%  '
$ws.Range("D20").Value = "'235.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.28%  '
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("D22").Value = "'2.160.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.48%  '
$ws.Range("D23").Value = "'1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.29%  '
$ws.Range("D24").Value = "'6.964"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.42%  '
$ws.Range("D25").Value = "'9.243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("D26").Value = "'165.71"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.15%  '
$ws.Range("D27").Value = "'18.68"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.41%  '
$ws.Range("D28").Value = "'0.1292"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.46%  '
$ws.Range("D29").Value = "'2.054"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.59%  '
$ws.Range("D30").Value = "'1.339"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.96%  '
$ws.Range("E31").Value = '  -1.91%  '
$ws.Range("D32").Value = "'4.291"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.26%  '
$ws.Range("D33").Value = "'4.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("D34").Value = "'0.05434"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.78%  '
$ws.Range("D35").Value = "'1.274"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.94%  '
$ws.Range("D36").Value = "'0.7328"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.14%  '
$ws.Range("D37").Value = "'2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.66%  '
$ws.Range("D38").Value = "'0.01933"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.55%  '
$ws.Range("D39").Value = "'2.760"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.39%  '
$ws.Range("D40").Value = "'6.225"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.99%  '
$ws.Range("D41").Value = "'0.4430"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").Value = "'72.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.45%  '
$ws.Range("D43").Value = "'1.929"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.65%  '
$ws.Range("D44").Value = "'1.002"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.11%  '
$ws.Range("D45").Value = "'0.8275"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = "'101.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = "'7.610"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("D48").Value = "'9.882"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").Value = "'2.063.33"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.19%  '
$ws.Range("D50").Value = "'36.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.18%  '
$ws.Range("D51").Value = "'0.05964"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.50%  '

Write-Host "Updated cryptos list"
